$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new article was published (row inserted at the top of the data, just
# below the header row). Insert a new row 2, which pushes every existing
# data row down by one (old row 2 -> row 3, old row 3 -> row 4, ...).
$ws.Rows.Item(2).Insert()

# The sheet's used range stays A1:E101 (100 data rows), so the last row
# (the old row 101, "What Are EMI Filters?") falls off the bottom and must
# be removed now that it lives at row 102.
$ws.Rows.Item(102).Delete()

# Populate the newly inserted row 2 with the new article's data.
$ws.Range("A2").Value = "Cryogenic Microwave Wafer-Scale Characterization Of Superconducting Resonators"
$ws.Range("B2").Value = "'7/9/2024"
$ws.Range("C2").Value = "In quantum computing, superconducting resonators are pivotal in enabling qubit readout and interaction. Explore their characterization at cryogenic temperatures using advanced wafer-scale measurement techniques."
$ws.Range("D2").Value = "quantum (ORG)"
$ws.Range("E2").Value = "https://www.rfglobalnet.com/doc/cryogenic-microwave-wafer-scale-characterization-of-superconducting-resonators-0001"
